$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that currently sits in paragraph 1
#    ("Rohit page").
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ------------------------------------------------------------------
# 2. Turn the single run "abcd" (2nd paragraph) into two runs:
#    "A" and "bcd". We first fix the capitalisation, then force a
#    run split by temporarily splitting the paragraph in two and
#    re-joining it (deleting the paragraph mark merges the two
#    paragraphs back into one without merging the runs back together).
# ------------------------------------------------------------------
$abcd = $d.Paragraphs.Item(2).Range
$aPos = $abcd.Start
$splitPos = $aPos + 1

# "a" -> "A"
$d.Range($aPos, $splitPos).Text = "A"

# split right after the "A"
$d.Range($splitPos, $splitPos).InsertParagraphAfter()

# merge the two paragraphs back together by deleting the paragraph
# mark that was just inserted - this keeps "A" and "bcd" as separate
# runs instead of re-merging them into a single run
$p2 = $d.Paragraphs.Item(2)
$markRange = $d.Range($p2.Range.End - 1, $p2.Range.End)
$markRange.Delete()

# ------------------------------------------------------------------
# 3. Add a brand-new paragraph ("Db mjb") right after the "Abcd"
#    paragraph.
# ------------------------------------------------------------------
$endOfP2 = $d.Paragraphs.Item(2).Range.End
$d.Range($endOfP2, $endOfP2).InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertAfter("Db mjb")

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark, collapsed, at the end of the
#    new "Db mjb" paragraph. Adding a bookmark with a collapsed range
#    sitting exactly in front of a paragraph mark mis-places it, so
#    we temporarily insert a placeholder character after the target
#    spot, anchor the bookmark there, and then remove the
#    placeholder again (the bookmark stays put).
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$endPos = $p3.Range.End - 1

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()
